$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.281.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.77%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.702.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.24%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'223.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.13%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.5301"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.51%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2659"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.90%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06583"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.33%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -4.78%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07621"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.57%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.497"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.06%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'1.936.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.27%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'WrappedEther"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'1.699.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.49%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.5786"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₅8142"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.74%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'67.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.01%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'27.272.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.88%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'215.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.24%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.06%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.610"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.77%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -3.40%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.968"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.28%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.14%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'144.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.51%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.703"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.61%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.97%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -3.02%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'16.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'0.05362"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.47%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.60%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.465"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.43%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.398"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.98%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.639"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.44%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.862"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.64%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.415"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.25%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.9455"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.51%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.5817"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.05%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01630"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.02%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.775"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.56%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.002"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'1.039.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.65%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.8395"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.77%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'100.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.83%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.844.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.22%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +1.15%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'57.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.07%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.4520"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.84%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.007"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.72%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'8.067"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.79%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.05226"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.81%  "
$ws.Range("E51").Style = "Normal"
